$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new header date cell BB1 (copy style from BA1)
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("BB1").Value = 45986

# Rows 2-71: BB column repeats the BA column value for that row
$ws.Range("BB2").Value = 0.1369474440744227
$ws.Range("BB3").Value = 1.366270496737897
$ws.Range("BB4").Value = -0.6985480789094254
$ws.Range("BB5").Value = -0.4148631161428114
$ws.Range("BB6").Value = -2.200693209579313
$ws.Range("BB7").Value = -4.008176425042492
$ws.Range("BB8").Value = 0.3143994836296855
$ws.Range("BB9").Value = 0.8129351329973105
$ws.Range("BB10").Value = 0.7286484790473651
$ws.Range("BB11").Value = 0.662286234074088
$ws.Range("BB12").Value = 2.210734389673945
$ws.Range("BB13").Value = 0.702250570321695
$ws.Range("BB14").Value = 0.6043683783303493
$ws.Range("BB15").Value = 1.531899235856926
$ws.Range("BB16").Value = 0.09102045989541807
$ws.Range("BB17").Value = 0.4274346081797518
$ws.Range("BB18").Value = 0.09961124419930911
$ws.Range("BB19").Value = 0.3110551286058296
$ws.Range("BB20").Value = 0.1249967678526218
$ws.Range("BB21").Value = 0.08642692979871924
$ws.Range("BB22").Value = -0.4125838437329037
$ws.Range("BB23").Value = -0.404671332649869
$ws.Range("BB24").Value = 0.793265767125348
$ws.Range("BB25").Value = 0.2975365817668774
$ws.Range("BB26").Value = 0.4497576285229741
$ws.Range("BB27").Value = 0.715037077548871
$ws.Range("BB28").Value = -0.05661231354093843
$ws.Range("BB29").Value = 0.1888053351092367
$ws.Range("BB30").Value = 0.6124584237519315
$ws.Range("BB31").Value = 0.170162512332567
$ws.Range("BB32").Value = 0.5212263459736306
$ws.Range("BB33").Value = 0.2407365673923465
$ws.Range("BB34").Value = 0.3602440716739608
$ws.Range("BB35").Value = 0.7156188677996056
$ws.Range("BB36").Value = 0.4660496629244335
$ws.Range("BB37").Value = 0.1637213432474738
$ws.Range("BB38").Value = 0.4177287092911968
$ws.Range("BB39").Value = 0.8997744569043959
$ws.Range("BB40").Value = 0.633132069676634
$ws.Range("BB41").Value = 0.735487593389081
$ws.Range("BB42").Value = 0.6069719124519395
$ws.Range("BB43").Value = 0.1320325676681762
$ws.Range("BB44").Value = 0.3861390137996494
$ws.Range("BB45").Value = -0.1032042178152466
$ws.Range("BB46").Value = 0.2066105200339621
$ws.Range("BB47").Value = 0.5
$ws.Range("BB48").Value = -0.2
$ws.Range("BB49").Value = 0.3
$ws.Range("BB50").Value = -0.1
$ws.Range("BB51").Value = -1.995361287679273
$ws.Range("BB52").Value = -9.697717272052344
$ws.Range("BB53").Value = 8.701161067295743
$ws.Range("BB54").Value = 0.5342924547835821
$ws.Range("BB55").Value = -1.684226516424943
$ws.Range("BB56").Value = 2.173700732922356
$ws.Range("BB57").Value = 1.669530332166502
$ws.Range("BB58").Value = -0.3471888372093019
$ws.Range("BB59").Value = 1.023580707979747
$ws.Range("BB60").Value = -0.1208922437305517
$ws.Range("BB61").Value = 0.4746827657805142
$ws.Range("BB62").Value = -0.5372032863913034
$ws.Range("BB63").Value = 0.2711818952007405
$ws.Range("BB64").Value = -0.07397500112315925
$ws.Range("BB65").Value = 0.1388009164387825
$ws.Range("BB66").Value = -0.4989801917293875
$ws.Range("BB67").Value = 0.2365428825421532
$ws.Range("BB68").Value = -0.2955890549112326
$ws.Range("BB69").Value = 0.1051437241507784
$ws.Range("BB70").Value = -0.2005382402049349
$ws.Range("BB71").Value = 0.3062442926496516

# Rows 72-82: BB column gets new forecast values (diverging from BA)
$ws.Range("BB72").Value = -0.2099029780610664
$ws.Range("BB73").Value = 0
$ws.Range("BB74").Value = 0.4020965382952741
$ws.Range("BB75").Value = 0.1648193638571329
$ws.Range("BB76").Value = 0.137888861916503
$ws.Range("BB77").Value = 0.2246336061890071
$ws.Range("BB78").Value = 0.1963670667155742
$ws.Range("BB79").Value = 0.1808293611245649
$ws.Range("BB80").Value = 0.1963958306960499
$ws.Range("BB81").Value = 0.1947212848581776
$ws.Range("BB82").Value = 0.1904909329952993

# Row 83: new row (copy style from A82 for the date cell)
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.1928178312549761
